$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits (order matters for shared-string table ordering) ---
# E5: "pozdějc - zhruba dneska" -> "pozdějc"
$ws.Range("E5").Value = "pozdějc"
# E7: "na pivko - něco mezi třema a šesti hodinama" -> "Napivkon"
$ws.Range("E7").Value = "Napivkon"
# E4: "dneska - takňák až to vyjde" -> "dnesk"
$ws.Range("E4").Value = "dnesk"

# --- Numeric edits (replace placeholder text with real numbers) ---
# F7: "Více než 180, méně než 360" -> 38.756
$ws.Range("F7").Value = 38.756
# F9: "Zadat poslední tři příchody" -> 15
$ws.Range("F9").Value = 15

# --- Column I cleanup ---
$ws.Range("I2").Value = "Jednotky smradu"
$ws.Range("I3").Value = "Jednotky prdu (zvuk)"
$ws.Range("I5").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("I7").Value = ""

# --- Restore automatic row heights on rows that previously had explicit heights ---
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(7).AutoFit()

# --- Selection / view state ---
$ws.Range("I4:I5").Select()
